$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 2.727484333333333
$ws.Range("H2").Value = 8.182453
$ws.Range("I2").Value = 0.03096049453772388
$ws.Range("J2").Value = 0.03096049453772388
$ws.Range("M2").Value = 547.450775
$ws.Range("N2").Value = 1642.352325
$ws.Range("O2").Value = 0.8253533007282613
$ws.Range("P2").Value = 0.8253533007282614
$ws.Range("Q2").Value = 1493.163412083692
$ws.Range("R2").Value = 13438.47070875323
$ws.Range("S2").Value = 0.02555334635888971
$ws.Range("T2").Value = 0.02555334635888971

# Row 3
$ws.Range("G3").Value = 2.727484333333333
$ws.Range("H3").Value = 8.182453
$ws.Range("I3").Value = 0.03096049453772388
$ws.Range("J3").Value = 0.03096049453772388
$ws.Range("O3").Value = 0.002183077622430991
$ws.Range("P3").Value = 0.002183077622430991
$ws.Range("Q3").Value = 3.949450046030445
$ws.Range("R3").Value = 35.54505041427401
$ws.Range("S3").Value = 0.00006758916280470192
$ws.Range("T3").Value = 0.00006758916280470195

# Row 4
$ws.Range("G4").Value = 2.727484333333333
$ws.Range("H4").Value = 8.182453
$ws.Range("I4").Value = 0.03096049453772388
$ws.Range("J4").Value = 0.03096049453772388
$ws.Range("M4").Value = 114.393852
$ws.Range("N4").Value = 343.181556
$ws.Range("O4").Value = 0.1724636216493076
$ws.Range("P4").Value = 0.1724636216493076
$ws.Range("Q4").Value = 312.007439159652
$ws.Range("R4").Value = 2808.066952436868
$ws.Range("S4").Value = 0.005339559016029465
$ws.Range("T4").Value = 0.005339559016029467

# Row 5
$ws.Range("I5").Value = 0.5986009007423507
$ws.Range("J5").Value = 0.5986009007423507
$ws.Range("M5").Value = 547.450775
$ws.Range("N5").Value = 1642.352325
$ws.Range("O5").Value = 0.8253533007282613
$ws.Range("P5").Value = 0.8253533007282614
$ws.Range("Q5").Value = 28869.33741771328
$ws.Range("R5").Value = 259824.0367594195
$ws.Range("S5").Value = 0.4940572292466095
$ws.Range("T5").Value = 0.4940572292466096

# Row 6
$ws.Range("I6").Value = 0.5986009007423507
$ws.Range("J6").Value = 0.5986009007423507
$ws.Range("O6").Value = 0.002183077622430991
$ws.Range("P6").Value = 0.002183077622430991
$ws.Range("R6").Value = 687.240288393692
$ws.Range("S6").Value = 0.001306792231177661
$ws.Range("T6").Value = 0.001306792231177661

# Row 7
$ws.Range("I7").Value = 0.5986009007423507
$ws.Range("J7").Value = 0.5986009007423507
$ws.Range("M7").Value = 114.393852
$ws.Range("N7").Value = 343.181556
$ws.Range("O7").Value = 0.1724636216493076
$ws.Range("P7").Value = 0.1724636216493076
$ws.Range("Q7").Value = 6032.459652468216
$ws.Range("R7").Value = 54292.13687221394
$ws.Range("S7").Value = 0.1032368792645635
$ws.Range("T7").Value = 0.1032368792645635

# Row 8
$ws.Range("G8").Value = 32.63402300000001
$ws.Range("H8").Value = 97.90206900000001
$ws.Range("I8").Value = 0.3704386047199253
$ws.Range("J8").Value = 0.3704386047199253
$ws.Range("M8").Value = 547.450775
$ws.Range("N8").Value = 1642.352325
$ws.Range("O8").Value = 0.8253533007282613
$ws.Range("P8").Value = 0.8253533007282614
$ws.Range("Q8").Value = 17865.52118271783
$ws.Range("R8").Value = 160789.6906444604
$ws.Range("S8").Value = 0.305742725122762
$ws.Range("T8").Value = 0.3057427251227621

# Row 9
$ws.Range("G9").Value = 32.63402300000001
$ws.Range("H9").Value = 97.90206900000001
$ws.Range("I9").Value = 0.3704386047199253
$ws.Range("J9").Value = 0.3704386047199253
$ws.Range("O9").Value = 0.002183077622430991
$ws.Range("P9").Value = 0.002183077622430991
$ws.Range("Q9").Value = 47.25469622844468
$ws.Range("R9").Value = 425.2922660560021
$ws.Range("S9").Value = 0.0008086962284486281
$ws.Range("T9").Value = 0.0008086962284486282

# Row 10
$ws.Range("G10").Value = 32.63402300000001
$ws.Range("H10").Value = 97.90206900000001
$ws.Range("I10").Value = 0.3704386047199253
$ws.Range("J10").Value = 0.3704386047199253
$ws.Range("M10").Value = 114.393852
$ws.Range("N10").Value = 343.181556
$ws.Range("O10").Value = 0.1724636216493076
$ws.Range("P10").Value = 0.1724636216493076
$ws.Range("Q10").Value = 3733.131597226597
$ws.Range("R10").Value = 33598.18437503937
$ws.Range("S10").Value = 0.0638871833687146
$ws.Range("T10").Value = 0.06388718336871461
